$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Message codes")
Write-Host "ActiveWindow:" $excel.ActiveWindow
$aw = $excel.ActiveWindow
Write-Host "ScrollRow before:" $aw.ScrollRow
$aw.ScrollRow = 63
Write-Host "ScrollRow after:" $aw.ScrollRow
